# Add API route for total profit (KPI)
# - New KPI row (row 29): Total profit = sum ((rate - cost - discount) * quantity)
# - C27 / C28 checks now compare against the live computed cell (B27 / B28)
#   instead of being hard-coded self-equal tautologies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the "Same result as with SQL" checks reference the actual computed
# totals instead of repeating the literal twice.
$ws.Range("C27").Formula = "=1782.04=B27"
$ws.Range("C28").Formula = "=856.09=B28"

# Insert a new row right after the "Total cost" row so row 28 (Total cost)
# keeps its position and a fresh row 29 is created below it.
$ws.Rows.Item(29).Insert()

# Seed the new row's formatting by copying row 28 (label / value / bool-check
# style) down, then overwrite the contents for the profit KPI.
$ws.Range("A28:C28").Copy()
$ws.Range("A29:C29").PasteSpecial()

$ws.Range("A29").Value = "Total profit = sum ((rate - cost - discount) * quantity)"
$ws.Range("B29").Formula = "=SUMPRODUCT((E2:E21 - J2:J21 - G2:G21) * D2:D21)"
$ws.Range("C29").Formula = "=583.91=B29"

# Move the active selection to C28 (matches the saved view state).
$ws.Range("C28").Select()
